$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K -> E:L), then copy the
# number formats/styles from the (now shifted) column E into the new
# column D so the new cells pick up the same styles as their row
# (date style for header rows, number style for data rows).
$ws.Columns("D:D").Insert()
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New (2018) column of data for the three statements (Income Statement,
# Balance Sheet, Cash Flow Statement).

# --- Income Statement (rows 7-35) ---
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 118900
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = -1400
$ws.Range("D17").Value2 = 29400
$ws.Range("D18").Value2 = 89500
$ws.Range("D20").Value2 = -46200
$ws.Range("D21").Value2 = 48000
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 43300
$ws.Range("D24").Value2 = 9600
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 33700
$ws.Range("D27").Value2 = 33700
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 46200
$ws.Range("D33").Value2 = 33700
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 33700

# --- Balance Sheet (rows 38-77) ---
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 6900
$ws.Range("D42").Value2 = 41200
$ws.Range("D43").Value2 = 0
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 7100
$ws.Range("D48").Value2 = 32400
$ws.Range("D49").Value2 = 110600
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 14100
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 2701300
$ws.Range("D57").Value2 = 0
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 0
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 56700
$ws.Range("D62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 2353000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 45000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 348300
$ws.Range("D77").Value2 = 0

# --- Cash Flow Statement (rows 80-102) ---
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 33700
$ws.Range("D83").Value2 = 4700
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 41400
$ws.Range("D91").Value2 = -2000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -93700
$ws.Range("D96").Value2 = -7700
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 55400
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 3100

$wb.Save()
